$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 15.12 = 61755.73 pesos`n✅ 61755.73 pesos = 15.0 = 968.91 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 66.149
$wsTasas.Range("O10").Value = 4085.08
$wsTasas.Range("N12").Value = 4117.43
